$wb = $excel.ActiveWorkbook

# Update the form_version value on the "settings" sheet (B3) from 1 to 20130408
$settings = $wb.Worksheets.Item("settings")
$settings.Range("B3").Value = 20130408

# Make "settings" the active sheet/tab (it becomes tabSelected in the saved file)
$settings.Activate()
$settings.Range("B3").Select()

$wb.Save()
